$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-21 Wednesday", 2)

$d.Content.Find.Execute("224÷9=24, 8", $true, $false, $false, $false, $false, $true, 1, $false, "293÷6=48, 5", 2)
$d.Content.Find.Execute("811÷9=90, 1", $true, $false, $false, $false, $false, $true, 1, $false, "563÷2=281, 1", 2)
$d.Content.Find.Execute("949÷7=135, 4", $true, $false, $false, $false, $false, $true, 1, $false, "795÷4=198, 3", 2)
$d.Content.Find.Execute("338÷7=48, 2", $true, $false, $false, $false, $false, $true, 1, $false, "237÷7=33, 6", 2)
$d.Content.Find.Execute("622÷8=77, 6", $true, $false, $false, $false, $false, $true, 1, $false, "697÷7=99, 4", 2)
$d.Content.Find.Execute("172÷6=28, 4", $true, $false, $false, $false, $false, $true, 1, $false, "207÷4=51, 3", 2)
$d.Content.Find.Execute("527÷9=58, 5", $true, $false, $false, $false, $false, $true, 1, $false, "801÷7=114, 3", 2)
$d.Content.Find.Execute("790÷6=131, 4", $true, $false, $false, $false, $false, $true, 1, $false, "689÷2=344, 1", 2)
$d.Content.Find.Execute("640÷4=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "713÷9=79, 2", 2)
$d.Content.Find.Execute("578÷4=144, 2", $true, $false, $false, $false, $false, $true, 1, $false, "134÷7=19, 1", 2)
$d.Content.Find.Execute("781÷4=195, 1", $true, $false, $false, $false, $false, $true, 1, $false, "276÷9=30, 6", 2)
$d.Content.Find.Execute("538÷8=67, 2", $true, $false, $false, $false, $false, $true, 1, $false, "749÷8=93, 5", 2)
$d.Content.Find.Execute("414÷8=51, 6", $true, $false, $false, $false, $false, $true, 1, $false, "852÷7=121, 5", 2)
$d.Content.Find.Execute("240÷2=120, 0", $true, $false, $false, $false, $false, $true, 1, $false, "428÷6=71, 2", 2)
$d.Content.Find.Execute("953÷8=119, 1", $true, $false, $false, $false, $false, $true, 1, $false, "267÷5=53, 2", 2)
$d.Content.Find.Execute("709÷3=236, 1", $true, $false, $false, $false, $false, $true, 1, $false, "381÷3=127, 0", 2)
$d.Content.Find.Execute("786÷3=262, 0", $true, $false, $false, $false, $false, $true, 1, $false, "112÷9=12, 4", 2)
$d.Content.Find.Execute("397÷7=56, 5", $true, $false, $false, $false, $false, $true, 1, $false, "942÷5=188, 2", 2)
$d.Content.Find.Execute("233÷3=77, 2", $true, $false, $false, $false, $false, $true, 1, $false, "148÷9=16, 4", 2)
$d.Content.Find.Execute("369÷9=41, 0", $true, $false, $false, $false, $false, $true, 1, $false, "291÷5=58, 1", 2)
$d.Content.Find.Execute("454÷2=227, 0", $true, $false, $false, $false, $false, $true, 1, $false, "688÷5=137, 3", 2)
$d.Content.Find.Execute("225÷8=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "492÷7=70, 2", 2)
$d.Content.Find.Execute("343÷5=68, 3", $true, $false, $false, $false, $false, $true, 1, $false, "490÷2=245, 0", 2)
$d.Content.Find.Execute("503÷3=167, 2", $true, $false, $false, $false, $false, $true, 1, $false, "257÷9=28, 5", 2)
$d.Content.Find.Execute("383÷4=95, 3", $true, $false, $false, $false, $false, $true, 1, $false, "941÷2=470, 1", 2)
